# Generate Report for Handback
# Applies the localization-status.xlsx handback update:
#   - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   - zh-cn / de-de sheets: fill in "Latest Target File" (I), "Latest Handback File" (J)
#     and "Latest Handback DateTime" (K) for both data rows, with a hyperlink on I.
#   - Column widths widened to fit the new (longer) text.

$wb = $excel.ActiveWorkbook

$hoTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/96dc900e071ca7fd2378845ce5e047c3ab8d8443/e2e/4d98daf5-3549-4f81-8d4b-c82bea531f69.md"
$hoDisplay = "4d98daf5-3549-4f81-8d4b-c82bea531f69.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text + widened zh-cn/de-de columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("I2").Value = $hoDisplay
$wsZhCn.Range("I3").Value = $hoDisplay
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $hoTarget, "", "", $hoDisplay)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $hoTarget, "", "", $hoDisplay)

$wsZhCn.Range("J2").Value = "4d98daf5-3549-4f81-8d4b-c82bea531f69.77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "4d98daf5-3549-4f81-8d4b-c82bea531f69.77613b4ad2a8b46aea42ac927a4e36deca9a5c06.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-23 19:06:53"
$wsZhCn.Range("K3").Value = "2016-08-23 19:06:53"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("I2").Value = $hoDisplay
$wsDeDe.Range("I3").Value = $hoDisplay
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $hoTarget, "", "", $hoDisplay)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $hoTarget, "", "", $hoDisplay)

$wsDeDe.Range("J2").Value = "4d98daf5-3549-4f81-8d4b-c82bea531f69.77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf"
$wsDeDe.Range("J3").Value = "4d98daf5-3549-4f81-8d4b-c82bea531f69.77613b4ad2a8b46aea42ac927a4e36deca9a5c06.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-23 19:07:02"
$wsDeDe.Range("K3").Value = "2016-08-23 19:07:02"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 39.1666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.1666666666667
